$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Last status check on: 28.02.2022 01:30"
$ws.Range("D5").Value = 0.4
$ws.Range("E5").Value = 44620.0521875
$ws.Range("E5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
